$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/value columns keep their text (string) formatting,
# same as the original "$x.xx" inline-string cells, rather than being
# auto-converted into currency numbers by Excel.
$ws.Range("C3:D6").NumberFormat = "@"

# Row 3: Iral74, 4, $3.40, $13.62 (was Idastidru52 data)
$ws.Range("A3").Value = "Iral74"
$ws.Range("C3").Value = "$3.40"
$ws.Range("D3").Value = "$13.62"

# Row 4: Idastidru52, 4, $3.86, $15.45 (was Iral74 data)
$ws.Range("A4").Value = "Idastidru52"
$ws.Range("C4").Value = "$3.86"
$ws.Range("D4").Value = "$15.45"

# Row 5: Aina42, 3, $3.07, $9.22
$ws.Range("A5").Value = "Aina42"
$ws.Range("C5").Value = "$3.07"
$ws.Range("D5").Value = "$9.22"

# Row 6: Aelin32, 3, $2.99, $8.98
$ws.Range("A6").Value = "Aelin32"
$ws.Range("C6").Value = "$2.99"
$ws.Range("D6").Value = "$8.98"
